# Remove the "ArrayList exercise" run from the subtitle placeholder on
# slide 1 (shape "Rectangle 3" / id 363523), while leaving the following
# line breaks and remaining text untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

$target = "ArrayList exercise"
$full = $tr.Text
$idx = $full.IndexOf($target)

if ($idx -ge 0) {
    $sub = $tr.Characters($idx + 1, $target.Length)
    $sub.Delete()
}
